# Applies proofing-mark (w:proofErr) annotations that Word's spell/grammar
# checker leaves behind after a pass over the document, plus adds a new
# line of text ("sfsfssfsdf") in what used to be an empty paragraph.
#
# The w:proofErr element has no surface in the Word object model (it is
# not reachable through any documented property/method), so the runs it
# wraps are rewritten in place via Range.InsertXML - this is exactly how
# Word itself persists the markers: the paragraph's full WordprocessingML
# is replaced with an equivalent copy that also carries the markers.

$d = $word.ActiveDocument

function Assert-ParaText($para, [string]$expected) {
    $actual = $para.Range.Text
    if ($actual -ne $expected) {
        throw "Unexpected paragraph text: expected [$expected] got [$actual]"
    }
}

function New-PkgXml([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# --- Paragraph 1: "DD" / "ssfsfs" / "asdfsffss" / "fsss" -------------------
# Wrap the whole run sequence with spellStart/spellEnd proofErr markers.
$p1 = $d.Paragraphs(1)
Assert-ParaText $p1 "DDssfsfsasdfsffssfsss`r"
$p1Xml = '<w:p w:rsidR="00AA579D" w:rsidRDefault="00955233">' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>DD</w:t></w:r>' +
    '<w:r w:rsidR="00270C57"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>ssfsfs</w:t></w:r>' +
    '<w:r w:rsidR="0025391F"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>asdfsffss</w:t></w:r>' +
    '<w:r w:rsidR="00AA579D"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>fsss</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
$p1.Range.InsertXML((New-PkgXml $p1Xml))

# --- Middle empty paragraph becomes "sfsfssfsdf" ----------------------------
# Of the 5 consecutive empty paragraphs after paragraph 1, the 3rd one
# (paragraph 4 overall) gains a new run of text, spell-checked and
# grammar-checked (spellStart/spellEnd + gramStart/gramEnd).
$p4 = $d.Paragraphs(4)
Assert-ParaText $p4 "`r"
$p4Xml = '<w:p>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>sfsfssfsdf</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '</w:p>'
$p4.Range.InsertXML((New-PkgXml $p4Xml))

# --- "sss" paragraph (with the _GoBack bookmark) ----------------------------
# Wrap the run + bookmark with spellStart/spellEnd + gramStart/gramEnd.
$p7 = $d.Paragraphs(7)
Assert-ParaText $p7 "sss`r"
$p7Xml = '<w:p w:rsidR="00AA579D" w:rsidRPr="00AA579D" w:rsidRDefault="00D4482D" w:rsidP="00AA579D">' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>sss</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '</w:p>'
$p7.Range.InsertXML((New-PkgXml $p7Xml))
